# Apply updated cryptocurrency price/volume data to Sheet1.
# Column D = Price, Column E = Volume(1h) change percentage.
# Only the cells whose values changed (per the source diff) are touched;
# all other cells (A, B, C columns, and unaffected D/E cells) are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.575.94"
$ws.Range("E2").Value = "  -1.61%  "
# Row 3
$ws.Range("D3").Value = "2.435.64"
$ws.Range("E3").Value = "  -2.27%  "
# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
# Row 5
$ws.Range("D5").Value = "'568.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.00%  "
# Row 6
$ws.Range("D6").Value = "'143.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.20%  "
# Row 7
$ws.Range("E7").Value = "  +0.10%  "
# Row 8
$ws.Range("D8").Value = "'0.531"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.08%  "
# Row 9
$ws.Range("D9").Value = "2.431.37"
$ws.Range("E9").Value = "  -2.38%  "
# Row 10
$ws.Range("E10").Value = "  -4.37%  "
# Row 11
$ws.Range("E11").Value = "  +1.30%  "
# Row 12
$ws.Range("D12").Value = "'5.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.14%  "
# Row 13
$ws.Range("D13").Value = "'0.354"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.00%  "
# Row 14
$ws.Range("D14").Value = "'26.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.20%  "
# Row 15
$ws.Range("E15").Value = "  -5.50%  "
# Row 16
$ws.Range("D16").Value = "2.877.25"
$ws.Range("E16").Value = "  -3.09%  "
# Row 17
$ws.Range("D17").Value = "62.483.85"
$ws.Range("E17").Value = "  -1.51%  "
# Row 18
$ws.Range("D18").Value = "2.441.82"
$ws.Range("E18").Value = "  -2.10%  "
# Row 19
$ws.Range("D19").Value = "'11.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.33%  "
# Row 20
$ws.Range("D20").Value = "'7.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.39%  "
# Row 21
$ws.Range("D21").Value = "'325.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.00%  "
# Row 22
$ws.Range("D22").Value = "'4.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.51%  "
# Row 23
$ws.Range("D23").Value = "'2.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.42%  "
# Row 25
$ws.Range("D25").Value = "'65.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.66%  "
# Row 26
$ws.Range("D26").Value = "'617.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.31%  "
# Row 27
$ws.Range("D27").Value = "'8.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.22%  "
# Row 28
$ws.Range("D28").Value = "0.0₃0989"
$ws.Range("E28").Value = "  -6.18%  "
# Row 29
$ws.Range("D29").Value = "2.564.87"
$ws.Range("E29").Value = "  -1.55%  "
# Row 30
$ws.Range("E30").Value = "  +0.64%  "
# Row 31
$ws.Range("E31").Value = "  -3.23%  "
# Row 32
$ws.Range("D32").Value = "'8.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.91%  "
# Row 33
$ws.Range("E33").Value = "  -2.06%  "
# Row 34
$ws.Range("D34").Value = "'0.137"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.71%  "
# Row 35
$ws.Range("D35").Value = "'5.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.32%  "
# Row 36
$ws.Range("D36").Value = "'1.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.23%  "
# Row 37
$ws.Range("E37").Value = "  +0.12%  "
# Row 38
$ws.Range("D38").Value = "'0.375"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.17%  "
# Row 39
$ws.Range("D39").Value = "'18.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.02%  "
# Row 40
$ws.Range("D40").Value = "'5.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.82%  "
# Row 41
$ws.Range("D41").Value = "'147.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.21%  "
# Row 42
$ws.Range("E42").Value = "  -5.44%  "
# Row 43
$ws.Range("E43").Value = "  -3.41%  "
# Row 44
$ws.Range("E44").Value = "  -0.02%  "
# Row 45
$ws.Range("E45").Value = "  +0.85%  "
# Row 46
$ws.Range("D46").Value = "'145.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.82%  "
# Row 47
$ws.Range("D47").Value = "'3.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.74%  "
# Row 48
$ws.Range("D48").Value = "'20.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.57%  "
# Row 49
$ws.Range("D49").Value = "'0.0528"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.30%  "
# Row 50
$ws.Range("D50").Value = "'0.596"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.66%  "
# Row 51
$ws.Range("D51").Value = "'0.0229"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.75%  "
